$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.887692451477051
$ws.Range("B1").Value = 5.216575622558594
$ws.Range("C1").Value = 6.87560510635376
$ws.Range("D1").Value = 10.51999282836914
$ws.Range("E1").Value = 5.36678409576416
